$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect to make edits, then restore protection afterward
$ws.Unprotect()

# Update the confidentiality / as-of-date disclaimer text in A7 (date 2021-03-24 -> 2021-03-25)
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Update Weight / Percent Change figures for EFA (row 2), EEM (row 3) and Total (row 4)
$ws.Range("D2").Value = 0.8448770091999219
$ws.Range("E2").Value = 0.007460698108180175

$ws.Range("D3").Value = 0.155122990800078
$ws.Range("E3").Value = 0.004643962848297267

$ws.Range("E4").Value = 0.007023757710375111

# Restore sheet protection
$ws.Protect()
